$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.636.87"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.821.85"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.82%  "
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4578"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3677"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07178"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8808"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07796"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.49"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.773.76"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.310"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.406"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "86.62"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.010"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008628"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.008"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.695.05"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.31"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.017"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.50"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.992"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.39"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.082"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.29"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.877"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08701"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.072"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.518"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7370"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.126"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.48%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.645"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.19%  "
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.006"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.083"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01946"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05126"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.905"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.033"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5056"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.209"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.009"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4650"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.02"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.50"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.598"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.73%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.49"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.84%  "
